# Add three new DMS test cases (DMS_SmartFolder, DMS_CreateView,
# DMS_TransmittalFolderView) to the DataFetchFlag / DataFetchXL config sheets,
# mirroring the existing rows, and update the active sheet/selection.

$wb = $excel.ActiveWorkbook
$wsFlag = $wb.Worksheets.Item("DataFetchFlag")
$wsXL = $wb.Worksheets.Item("DataFetchXL")

$testNames = @("DMS_SmartFolder", "DMS_CreateView", "DMS_TransmittalFolderView")
$targetRows = @(4, 5, 6)

# Same values already used for the existing DmsDocumentLibrary row (B3).
$dmsPathValue = $wsXL.Range("B3").Value2
$dmsHyperlinkAddress = "file:///\\src\com\proj\suiteDOCS\testdata\DmsDocumentLibrary.xlsx"

for ($i = 0; $i -lt $testNames.Count; $i++) {
    $row = $targetRows[$i]
    $name = $testNames[$i]

    # DataFetchFlag sheet: TestCaseName / DataFetchFlag ("XL")
    $wsFlag.Range("A$row").Value = $name
    $wsFlag.Range("B$row").Value = "XL"

    # DataFetchXL sheet: TestCaseName / ExcelDataSheetPath / FirstSheetName
    $wsXL.Range("A$row").Value = $name
    $wsXL.Range("B$row").Value = $dmsPathValue
    $wsXL.Range("C$row").Value = $name

    # Hyperlink the path cell to the DmsDocumentLibrary workbook, same as B2/B3.
    $wsXL.Hyperlinks.Add($wsXL.Range("B$row"), $dmsHyperlinkAddress)
    $wsXL.Range("B$row").Style = "Hyperlink"
}

# Update selection/active sheet: DataFetchXL becomes inactive (selection A6),
# DataFetchFlag becomes the active tab (selection B6).
$wsXL.Range("A6").Select()
$wsFlag.Range("B6").Select()
